$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Set ADC to 12/10/8/6-bit" commands are being collapsed into a single
# "Set ADC resolution" opcode, and "Set RTC time"/"Set RTC date" become their
# own distinct commands instead of duplicating "Get RTC time"/"Get RTC date".
# Everything below shifts up into rows 27-37 and the now-superfluous trailing
# row 38 is cleared.

# Column A (command names), top to bottom.
$ws.Range("A27").Value = "Set ADC resolution"
$ws.Range("A28").Value = "Send over USB"
$ws.Range("A29").Value = "Set RTC time"
$ws.Range("A30").Value = "Set RTC date"
$ws.Range("A31").Value = "Get RTC time"
$ws.Range("A32").Value = "Get RTC date"
$ws.Range("A33").Value = "Get DAC resolution"
$ws.Range("A34").Value = "Set DAC to 12-bit"
$ws.Range("A35").Value = "Set DAC to 8-bit"
$ws.Range("A36").Value = "Output on DAC 0"
$ws.Range("A37").Value = "Output on DAC 1"

# Column D (data width), top to bottom.
$ws.Range("D27").Value = "2-bit"
$ws.Range("D28").Value = "8-bit   ?"
$ws.Range("D29").Value = "24-bit"
$ws.Range("D30").Value = "24-bit"
$ws.Range("D31").Value = "24-bit"
$ws.Range("D32").Value = "24-bit"
$ws.Range("D33").Value = "1-bit"
$ws.Range("D34").Value = "1-bit"
$ws.Range("D35").Value = "1-bit"
$ws.Range("D36").Value = "DAC resolution"
$ws.Range("D37").Value = "DAC resolution"

# Column B (opcode numbers) renumbered now that four old rows collapsed to three.
$ws.Range("B29").Value = 28
$ws.Range("B30").Value = 29
$ws.Range("B31").Value = 30
$ws.Range("B32").Value = 31
$ws.Range("B33").Value = 32
$ws.Range("B34").Value = 33
$ws.Range("B35").Value = 34
$ws.Range("B36").Value = 35
$ws.Range("B37").Value = 36

# Old row 38 ("Output on DAC 1" duplicate) is now redundant - clear it so the
# row element disappears from the sheet XML.
$ws.Range("A38:D38").ClearContents()

# Update the view state to match the saved selection.
$ws.Range("D21").Select()
